$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 12:22"

# Row 4 (Madrid) - updated totals
$ws.Range("B4").Value = 42450
$ws.Range("C4").Value = 19836
$ws.Range("D4").Value = 17028
$ws.Range("E4").Value = 5586

# Row 9 (Ciudad Real) - updated totals
$ws.Range("B9").Value = 3949
$ws.Range("C9").Value = 733
$ws.Range("D9").Value = 2873
$ws.Range("E9").Value = 343

# Row 10 (Valencia/Valencia) - updated totals
$ws.Range("B10").Value = 3467
$ws.Range("C10").Value = 450
$ws.Range("D10").Value = 2811
$ws.Range("E10").Value = 206

# Row 11 (Navarra) - updated totals
$ws.Range("B11").Value = 2951
$ws.Range("C11").Value = 1061
$ws.Range("D11").Value = 1713
$ws.Range("E11").Value = 177

# Rows 13-15: province order changed (Alacant/Alicante moved earlier) plus updated totals
$ws.Range("A13").Value = "Alacant/Alicante"
$ws.Range("B13").Value = 2803
$ws.Range("C13").Value = 469
$ws.Range("D13").Value = 2042
$ws.Range("E13").Value = 292

$ws.Range("A14").Value = "Castilla-La Mancha"
$ws.Range("B14").Value = 2780
$ws.Range("C14").Value = 71
$ws.Range("D14").Value = 2446
$ws.Range("E14").Value = 263

$ws.Range("A15").Value = "Araba/Alava"
$ws.Range("B15").Value = 2744
$ws.Range("C15").Value = 3728
$ws.Range("D15").Value = 4707
$ws.Range("E15").Value = 212

# Rows 24-27: province order changed (Cantabria moved earlier) plus updated totals
$ws.Range("A24").Value = "Cantabria"
$ws.Range("B24").Value = 1572
$ws.Range("C24").Value = 175
$ws.Range("D24").Value = 1305
$ws.Range("E24").Value = 92

$ws.Range("A25").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B25").Value = 1560
$ws.Range("C25").Value = 3728
$ws.Range("D25").Value = 4707
$ws.Range("E25").Value = 82

$ws.Range("A26").Value = "Granada"
$ws.Range("B26").Value = 1550
$ws.Range("C26").Value = 182
$ws.Range("D26").Value = 1240
$ws.Range("E26").Value = 128

$ws.Range("A27").Value = "Pontevedra"
$ws.Range("B27").Value = 1536
$ws.Range("C27").Value = 333
$ws.Range("D27").Value = 1411
$ws.Range("E27").Value = 30

# Rows 37-39: province order changed (Castello/Castellon moved earlier) plus updated totals
$ws.Range("A37").Value = "Castello/Castellon"
$ws.Range("B37").Value = 899
$ws.Range("C37").Value = 142
$ws.Range("D37").Value = 668
$ws.Range("E37").Value = 89

$ws.Range("A38").Value = "Guadalajara"
$ws.Range("B38").Value = 897
$ws.Range("C38").Value = 1353
$ws.Range("D38").Value = 8547
$ws.Range("E38").Value = 128

$ws.Range("A39").Value = "Cadiz"
$ws.Range("B39").Value = 881
$ws.Range("C39").Value = 109
$ws.Range("D39").Value = 740
$ws.Range("E39").Value = 32
